# All Export .xls File Names Changed to .xlsx
# (the actual workbook edit: remove the obsolete "EndTime" column from the
# Users export template, shifting the trailing "don't remove" helper column
# left from M to L)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Column L currently holds the header "EndTime" (L1) and an example value
# (L2) that are no longer used by the export. Deleting the entire column
# shifts the following "don't remove" helper column (M) left to become the
# new column L, exactly matching the target workbook layout.
$ws.Range("L:L").Delete()

# Restore selection on the (now shifted) helper column, matching how Excel
# leaves the selection positioned on the column after a delete operation.
$ws.Range("L1:L1048576").Select()
